# Insert a new price-record row at row 125 (Femacal de La Calera - Espinaca),
# pushing the existing rows 125-185 down to 126-186.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(125).Insert()

# Populate the newly inserted row with the new record's data.
$ws.Range("A125").Value = 3
$ws.Range("B125").Value = "Femacal de La Calera"
$ws.Range("C125").Value = "Coquimbo"
$ws.Range("D125").Value = 44460
$ws.Range("E125").Value = 5
$ws.Range("F125").Value = 100112012
$ws.Range("G125").Value = "Espinaca"
$ws.Range("H125").Value = "Sin especificar"
$ws.Range("I125").Value = "Primera"
$ws.Range("J125").Value = 130
$ws.Range("K125").Value = 2500
$ws.Range("L125").Value = 3000
$ws.Range("M125").Value = 2769
$ws.Range("N125").Value = "$/docena de atados (3 kilos)"
$ws.Range("O125").Value = "Provincia de Quillota"
$ws.Range("P125").Value = 923
$ws.Range("Q125").Value = 3
$ws.Range("R125").Value = "Hortaliza"
